$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2:K41").ClearContents()
